$wb = $excel.ActiveWorkbook

# --- Cloud 01 sheet ---
$ws1 = $wb.Worksheets.Item("Cloud 01")
$ws1.Range("B3").Value = "password12"
$ws1.Range("B4").Value = "APK File/IvyDMS_PNGIndia_1575_9922.apk"
$ws1.Range("B13").Value = "TNPKL-23-I000517"
$ws1.Range("A6").Copy() | Out-Null
$ws1.Range("D6").PasteSpecial(-4122) | Out-Null

$ws1.Range("B3").Select()

# --- Cloud 02 sheet ---
$ws2 = $wb.Worksheets.Item("Cloud 02")
$ws2.Range("B1").Value = "https://cloud02-in.ivydms.com/web/DMS"
$ws2.Range("B3").Value = "password12"
$ws2.Range("B3").Select()

# --- Cloud03 sheet ---
$ws3 = $wb.Worksheets.Item("Cloud03")
$ws3.Range("B1").Value = "https://cloud03-in.ivydms.com/web/DMS"
$ws3.Range("B3").Value = "password12"
$ws3.Range("B3").Select()
